$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Executed" column (D) values for the testcases
$ws.Range("D2").Value = "no"
$ws.Range("D3").Value = "no"
$ws.Range("D5").Value = "no"
$ws.Range("D6").Value = "yes"

# Move the active cell selection to D6
$ws.Activate()
$ws.Range("D6").Select()
